{"js": "// Update the GPA shown on the resume from \"3.85 GPA\" to \"3.90 GPA\".\n// The paragraph text (across runs) reads:\n//   \"Arizona State University, Tempe, AZ\" <tab> \"3.\" \"8\" \"5 GPA\"\n// We locate the single paragraph containing \"GPA\", then within that\n// paragraph replace the \"8\" run with \"90\" and the \"5 GPA\" run with\n// \" GPA\" (keeping the leading space) so the visible text becomes\n// \"3.90 GPA\" while leaving every other run/paragraph untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet gpaParagraph = null;\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(\"GPA\") !== -1) {\n    gpaParagraph = p;\n    break;\n  }\n}\n\nif (!gpaParagraph) {\n  throw new Error(\"Could not find the paragraph containing 'GPA'.\");\n}\n\n// Scope the searches to just this paragraph so they are unambiguous.\nconst digitResults = gpaParagraph.search(\"8\", { matchWholeWord: false, matchCase: true });\nconst tailResults = gpaParagraph.search(\"5 GPA\", { matchWholeWord: false, matchCase: true });\ndigitResults.load(\"items\");\ntailResults.load(\"items\");\nawait context.sync();\n\nif (digitResults.items.length !== 1 || tailResults.items.length !== 1) {\n  throw new Error(\"Unexpected number of GPA text matches found.\");\n}\n\n// Replace \"5 GPA\" -> \" GPA\" first, then \"8\" -> \"90\" (order doesn't matter\n// for correctness here since the ranges don't overlap, but doing the\n// later-in-text one first keeps offsets stable either way).\ntailResults.items[0].insertText(\" GPA\", Word.InsertLocation.replace);\ndigitResults.items[0].insertText(\"90\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Update the GPA shown on the resume from \"3.85 GPA\" to \"3.90 GPA\".\n# The paragraph text (spanning multiple runs) reads:\n#   \"Arizona State University, Tempe, AZ\" <tab> \"3.\" \"8\" \"5 GPA\"\n# Locate that paragraph, then, scoped to just its Range, find-and-replace\n# the \"8\" with \"90\" and \"5 GPA\" with \" GPA\" (keeping the leading space)\n# so the visible text becomes \"3.90 GPA\" without touching any other\n# paragraph/run in the document.\n\n$d = $word.ActiveDocument\n\n$gpaParagraph = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*GPA*\") {\n        $gpaParagraph = $p\n        break\n    }\n}\n\nif ($null -eq $gpaParagraph) {\n    throw \"Could not find the paragraph containing 'GPA'.\"\n}\n\n$paraStart = $gpaParagraph.Range.Start\n$paraEnd = $gpaParagraph.Range.End\n\n# Replace \"5 GPA\" -> \" GPA\" (the later text), then \"8\" -> \"90\".\n$tailRange = $d.Range($paraStart, $paraEnd)\n$tailFind = $tailRange.Find\n$tailFind.ClearFormatting()\n$tailFind.Replacement.ClearFormatting()\n$tailFind.Text = \"5 GPA\"\n$tailFind.Replacement.Text = \" GPA\"\n$tailFound = $tailFind.Execute([ref]\"5 GPA\", $false, $false, $false, $false, $false, $true, 1, $false, [ref]\" GPA\", 1)\nif (-not $tailFound) {\n    throw \"Could not find '5 GPA' text to replace.\"\n}\n\n$digitRange = $d.Range($paraStart, $paraEnd)\n$digitFind = $digitRange.Find\n$digitFind.ClearFormatting()\n$digitFind.Replacement.ClearFormatting()\n$digitFind.Text = \"8\"\n$digitFind.Replacement.Text = \"90\"\n$digitFound = $digitFind.Execute([ref]\"8\", $false, $false, $false, $false, $false, $true, 1, $false, [ref]\"90\", 1)\nif (-not $digitFound) {\n    throw \"Could not find the '8' digit text to replace.\"\n}\n"}
